# Apply the "Fixed update to excel issue" changes:
#  1. Rename header B1 on "Weekly Quantity" from "Requested quantity" to "Weekly_PO_Qty"
#  2. Rename header B1 on "Monthly Trend" from "Requested quantity" to "Monthly_PO_Qty"
#  3. Add a new "PO Forecast" worksheet (after "Monthly Trend") with forecast data

$wb = $excel.ActiveWorkbook

# --- 1. Weekly Quantity header rename -------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- 2. Monthly Trend header rename ----------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3. New "PO Forecast" sheet --------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows
$data = @(
    @(45067.99999999999, 40,  40.00158742273024,  40.00158742287264),
    @(45102.99999999999, 200, 200.0015346228813,  200.0015346230284),
    @(45109.99999999999, 232, 232.0015240536762,  232.0015240707613),
    @(45116.99999999999, 264, 264.0015132286008,  264.001513718058),
    @(45123.99999999999, 296, 296.0015023500455,  296.0015034817582),
    @(45130.99999999999, 328, 328.0014913410088,  328.0014932995331),
    @(45137.99999999999, 360, 360.0014803091519,  360.001483274011),
    @(45144.99999999999, 392, 392.0014690877254,  392.0014733266294),
    @(45151.99999999999, 424, 424.0014577286921,  424.0014634709423),
    @(45158.99999999999, 456, 456.001446271056,   456.0014535852833)
)

$row = 2
foreach ($rec in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $rec[0]
    $wsForecast.Cells.Item($row, 2).Value = $rec[1]
    $wsForecast.Cells.Item($row, 3).Value = $rec[2]
    $wsForecast.Cells.Item($row, 4).Value = $rec[3]
    $row++
}

# Date formatting for column A (rows 2-11), matching the other sheets' date style
$wsForecast.Range("A2:A11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
